$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a value to be stored as literal text (not auto-converted to a
# number by Excel), while leaving the cell's style index unchanged afterward.
function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "60.015.60"
$ws.Range("E2").Value = "  -0.96%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.418.44"
$ws.Range("E3").Value = "  -1.26%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.01%  "

# Row 5 - BNB
Set-TextValue "D5" "552.27"
$ws.Range("E5").Value = "  -0.94%  "

# Row 6 - Solana
Set-TextValue "D6" "137.24"
$ws.Range("E6").Value = "  -1.50%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.03%  "

# Row 8 - XRP
Set-TextValue "D8" "0.595"
$ws.Range("E8").Value = "  +3.96%  "

# Row 9 - Dogecoin
Set-TextValue "D9" "0.105"
$ws.Range("E9").Value = "  -1.71%  "

# Row 10 - Toncoin
$ws.Range("E10").Value = "  -2.42%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  -0.95%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  -2.31%  "

# Row 13 - Avalanche
Set-TextValue "D13" "25.21"
$ws.Range("E13").Value = "  +0.67%  "

# Row 14 - WrappedliquidstakedEther2.0
Set-TextValue "D14" "2.850.07"
$ws.Range("E14").Value = "  -1.15%  "

# Row 15 - WrappedBTC
Set-TextValue "D15" "59.966.96"
$ws.Range("E15").Value = "  -0.84%  "

# Row 16 - ShibaInu
$ws.Range("E16").Value = "  -1.92%  "

# Row 17 - WrappedEther
Set-TextValue "D17" "2.415.96"
$ws.Range("E17").Value = "  -0.44%  "

# Row 18 - Chainlink
Set-TextValue "D18" "11.29"
$ws.Range("E18").Value = "  -1.65%  "

# Row 19 - Polkadot
$ws.Range("E19").Value = "  -0.33%  "

# Row 20 - BitcoinCash
Set-TextValue "D20" "327.22"
$ws.Range("E20").Value = "  -2.78%  "

# Row 21 - Uniswap
$ws.Range("E21").Value = "  -3.63%  "

# Row 22 - Dai
$ws.Range("E22").Value = "  +0.01%  "

# Row 24 - Kaspa
$ws.Range("E24").Value = "  +2.74%  "

# Row 25 - InternetComputer(DFINITY)
Set-TextValue "D25" "8.60"
$ws.Range("E25").Value = "  +0.45%  "

# Row 26 - Binance-PegBSC-USD
$ws.Range("E26").Value = "  +0.09%  "

# Row 27 - Fetch.AI
$ws.Range("E27").Value = "  +1.02%  "

# Row 28 - PEPE
Set-TextValue "D28" "0.0₃0778"
$ws.Range("E28").Value = "  -2.72%  "

# Row 29 - PancakeSwap
$ws.Range("E29").Value = "  -2.43%  "

# Row 30 - Monero
Set-TextValue "D30" "168.24"
$ws.Range("E30").Value = "  -1.45%  "

# Row 31 - Aptos
Set-TextValue "D31" "6.06"
$ws.Range("E31").Value = "  -4.11%  "

# Row 32 - SuiNetwork
Set-TextValue "D32" "1.05"
$ws.Range("E32").Value = "  +2.14%  "

# Row 33 - EthereumClassic
$ws.Range("E33").Value = "  -1.50%  "

# Row 34 - USDe
$ws.Range("E34").Value = "  -0.02%  "

# Row 35 - ImmutableX
Set-TextValue "D35" "1.31"
$ws.Range("E35").Value = "  -0.99%  "

# Row 36 - FirstDigitalUSD
$ws.Range("E36").Value = "  -0.02%  "

# Row 37 - NEARProtocol
$ws.Range("E37").Value = "  -2.44%  "

# Row 38 - was Stacks, now Bittensor
$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D38" "325.43"
$ws.Range("E38").Value = "  +1.98%  "

# Row 39 - was Bittensor, now Stacks
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D39" "1.61"
$ws.Range("E39").Value = "  -2.50%  "

# Row 40 - PolygonEcosystemToken
Set-TextValue "D40" "0.405"
$ws.Range("E40").Value = "  -3.19%  "

# Row 41 - Filecoin
$ws.Range("E41").Value = "  -1.94%  "

# Row 42 - Aave
Set-TextValue "D42" "140.34"
$ws.Range("E42").Value = "  -2.94%  "

# Row 43 - Stellar
Set-TextValue "D43" "0.0972"
$ws.Range("E43").Value = "  +0.64%  "

# Row 44 - InjectiveProtocol
Set-TextValue "D44" "19.62"
$ws.Range("E44").Value = "  -2.07%  "

# Row 45 - Hedera
Set-TextValue "D45" "0.0516"
$ws.Range("E45").Value = "  -1.94%  "

# Row 46 - Mantle
Set-TextValue "D46" "0.578"
$ws.Range("E46").Value = "  +0.47%  "

# Row 47 - Polygon
Set-TextValue "D47" "0.399"
$ws.Range("E47").Value = "  -2.28%  "

# Row 48 - VeChain
$ws.Range("E48").Value = "  -2.07%  "

# Row 49 - WhiteBITCoin
$ws.Range("E49").Value = "  -0.07%  "

# Row 50 - dogwifhat
$ws.Range("E50").Value = "  -4.72%  "

# Row 51 - ZEEBU
$ws.Range("E51").Value = "  -1.32%  "
